# Applies the commit: appends 11 new sample data rows (1003-1013) below the
# existing data range on Sheet1, re-using the existing "asdf" shared string
# and introducing a new shared string "qwer", written into both column A
# and column B of each new row. This mirrors the author's pattern of rows
# already present (e.g. row 2, which holds "asdf" in both A and B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("qwer", "qwer", "asdf", "asdf", "asdf", "qwer", "qwer", "asdf", "asdf", "asdf", "asdf")

$startRow = 1003
for ($i = 0; $i -lt $values.Length; $i++) {
    $r = $startRow + $i
    $val = $values[$i]
    $ws.Cells.Item($r, 1).Value = $val
    $ws.Cells.Item($r, 2).Value = $val
}
